$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: Id changes from 7238787 to 92649821
$ws.Range("A2").Value = 92649821

# N2 (Metod): new empty cell (kept blank but present)
$ws.Range("N2").NumberFormat = "@"

# P2 (Lokalnamn): text changed
$ws.Range("P2").Value = "Kärrgårda,  700 m NV , Sk"

# S2 (Noggrannhet): 25 -> 10
$ws.Range("S2").Value = 10

# X2 (Externid): new value
$ws.Range("X2").Value = "M-Äng-0128"

# AC2 (Publik kommentar): text changed
$ws.Range("AC2").Value = "En gammal lokal som ansågs utdöd, men jag lyckades hitta ett tynande ex i skogen. Hårt trängt exempla. Behöver röjas om arten skall överleva"

# AF2 (Bestämningsmetod): new empty cell
$ws.Range("AF2").NumberFormat = "@"

# AW2 (Rapportör): Ulf Ryde -> Charlotte Wigermo
$ws.Range("AW2").Value = "Charlotte Wigermo"

# AX2 (Observatörer): stays Ulf Ryde (already correct, set explicitly for safety)
$ws.Range("AX2").Value = "Ulf Ryde"

# AY2 (Projektnamn): empty -> Floraväkteri Sverige
$ws.Range("AY2").Value = "Floraväkteri Sverige"
